# Excel_plan.xlsx update
# - Roadmap (row 20) items re-ordered / renamed:
#     "Adatbázis tervezés" step removed, "Adatbázis létrehozás" promoted
#     new "C# Osztályok létrehozása" step added
#     "Adatbevitel az Sql adatbázisba C#-al" renamed -> "Adat írás/olvasás az Sql adatbázisba C#-al"
#     "Dokumentáció" / "Extrák" steps shifted right (new K column)
# - Assorted empty placeholder/formatting-only cells removed (rows 6 & 8, stray column F cells, etc.)
# - Header/legend highlight colours swapped from the old theme tints to flat RGB fills
# - Selection moved to G21

function RGBColor($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the now-empty formatting-only cells (no value, fill reset to none)
# ---------------------------------------------------------------------------
$ws.Range("F1").Clear()
$ws.Range("E2:F2").Clear()
$ws.Range("D3:F3").Clear()
$ws.Range("E4:F4").Clear()
$ws.Range("F5").Clear()
$ws.Range("A6:F6").Clear()
$ws.Range("E7:F7").Clear()
$ws.Range("A8:F8").Clear()
$ws.Range("D9:F9").Clear()
$ws.Range("E10:F10").Clear()
$ws.Range("F11").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("A15").Clear()
$ws.Range("A16:B16").Clear()
$ws.Range("A17:B17").Clear()
$ws.Range("D17").Clear()
$ws.Range("A18:B18").Clear()
$ws.Range("D18").Clear()

# ---------------------------------------------------------------------------
# 2) Roadmap row (row 20) content updates
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = "Adatbázis létrehozás"
$ws.Range("C20").Value = "C# Osztályok létrehozása"
$ws.Range("D20").Value = "C# Összekötése az Sql adatbázissal"
$ws.Range("E20").Value = "C# Wpf alap UI"
$ws.Range("F20").Value = "C# Backend, adatkezelés"
$ws.Range("G20").Value = "Adat írás/olvasás az Sql adatbázisba C#-al"
$ws.Range("H20").Value = "Adat titkosítás kivitelezése"
$ws.Range("I20").Value = "App UI és design"
$ws.Range("J20").Value = "Dokumentáció"
$ws.Range("K20").Value = "Extrák"

# New K20 cell needs the same look as the other roadmap "chip" cells (bold, red fill)
$ws.Range("K20").Font.Bold = $true
$ws.Range("K20").Interior.Color = RGBColor 255 80 80

# ---------------------------------------------------------------------------
# 3) Re-colour the highlighted ranges (theme tints -> flat RGB colours)
# ---------------------------------------------------------------------------
$green = RGBColor 0 176 80
$yellow = RGBColor 255 255 0

# ER-diagram / Sql adatbázis column headers
$ws.Range("B1:C1").Interior.Color = $green

# Roadmap "Adatbázis létrehozás" / "C# Osztályok létrehozása" / "C# Összekötése az Sql adatbázissal"
$ws.Range("B20:D20").Interior.Color = $green

# Sql táblák column legends
$ws.Range("B13:D15").Interior.Color = $yellow
$ws.Range("C16:D16").Interior.Color = $yellow
$ws.Range("C17").Interior.Color = $yellow
$ws.Range("C18").Interior.Color = $yellow

# ---------------------------------------------------------------------------
# 4) Selection
# ---------------------------------------------------------------------------
$ws.Range("G21").Select()
